# Update "想去人数" (number of people interested) values for a handful of
# events. The same events are listed on both the "展览" sheet and the
# "全部类型" sheet, so both need to be updated in sync.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Row = 8;  Value = 19 },
    @{ Row = 11; Value = 569 },
    @{ Row = 13; Value = 13477 },
    @{ Row = 17; Value = 5550 },
    @{ Row = 18; Value = 5580 },
    @{ Row = 19; Value = 55 }
)

$ws1 = $wb.Worksheets.Item("展览")
foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 6).Value = $u.Value
}

$updates4 = @(
    @{ Row = 30; Value = 19 },
    @{ Row = 33; Value = 569 },
    @{ Row = 35; Value = 13477 },
    @{ Row = 40; Value = 5550 },
    @{ Row = 41; Value = 5580 },
    @{ Row = 42; Value = 55 }
)

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($u in $updates4) {
    $ws4.Cells.Item($u.Row, 6).Value = $u.Value
}
